# Apply the "DOSE_INFO test case working and passing" edit to the
# dose_info_scanned.xlsx fixture's "Files" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Files")

# Row 2 (image-00000.dcm scan): DICOM:StudyDate (O2) was mis-recorded as
# 20200312 -- correct it to the real study date 20161223, matching the
# Session value in G2. Both become text (quote-prefixed) rather than a
# bare number, so Excel doesn't reformat/re-interpret them as dates.
$ws.Range("G2").Value = "'20161223"
$ws.Range("O2").Value = "'20161223"

# Row 3 (image-00001.dcm scan): same fix.
$ws.Range("G3").Value = "'20161223"
$ws.Range("O3").Value = "'20161223"

# Row 4: this fixture row now documents the dose_info.dcm file (previously
# was referencing an unrelated img-with-pdf.dcm / encapsulated-report
# example) -- update filename and status message to reflect the DOSE_INFO
# test case.
$ws.Range("C4").Value = "dose_info.dcm"
$ws.Range("E4").Value = "DICOM has ImageType DOSE_INFO"

# Move the active selection to E1 (matches the saved view state).
$ws.Range("E1").Select()
